$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '24.054.14'
$ws.Range('E2').Value = '  -4.21%  '
$ws.Range('D3').Value = '1.639.94'
$ws.Range('E3').Value = '  -4.04%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.002'
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '307.19'
$ws.Range('E5').Value = '  -3.35%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '1.002'
$ws.Range('E6').Value = '  +0.09%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.3906'
$ws.Range('E7').Value = '  -2.64%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.3835'
$ws.Range('E8').Value = '  -5.13%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '1.003'
$ws.Range('E9').Value = '  -0.11%  '
$ws.Range('B10').Value = 'OKB'
$ws.Range('C10').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '49.03'
$ws.Range('E10').Value = '  -7.31%  '
$ws.Range('B11').Value = 'Polygon'
$ws.Range('C11').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '1.336'
$ws.Range('E11').Value = '  -9.30%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.08396'
$ws.Range('E12').Value = '  -4.92%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '23.61'
$ws.Range('E13').Value = '  -9.42%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '7.063'
$ws.Range('E14').Value = '  -5.84%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.00001273'
$ws.Range('E15').Value = '  -6.15%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '7.431'
$ws.Range('E16').Value = '  -6.91%  '
$ws.Range('D17').Value = '1.639.18'
$ws.Range('E17').Value = '  -3.72%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '94.83'
$ws.Range('E18').Value = '  -1.55%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.06890'
$ws.Range('E19').Value = '  -4.38%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '21.02'
$ws.Range('E20').Value = '  +0.85%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '6.890'
$ws.Range('E21').Value = '  -5.81%  '
$ws.Range('E22').Value = '  -0.01%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '13.55'
$ws.Range('E23').Value = '  -5.65%  '
$ws.Range('D24').Value = '24.051.39'
$ws.Range('E24').Value = '  -4.17%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.321'
$ws.Range('E25').Value = '  -3.32%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '2.690'
$ws.Range('E26').Value = '  -8.90%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '22.24'
$ws.Range('E27').Value = '  -5.83%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '156.92'
$ws.Range('E28').Value = '  -3.85%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '8.699'
$ws.Range('E29').Value = '  +4.15%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '140.32'
$ws.Range('E30').Value = '  -7.80%  '
$ws.Range('E31').Value = '  -14.98%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '2.442'
$ws.Range('E32').Value = '  -9.40%  '
$ws.Range('D33').Value = '1.818.61'
$ws.Range('E33').Value = '  -3.91%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '6.965'
$ws.Range('E34').Value = '  -3.48%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.07954'
$ws.Range('E35').Value = '  -7.75%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.02885'
$ws.Range('E36').Value = '  -9.25%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.9499'
$ws.Range('E37').Value = '  -9.39%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.2668'
$ws.Range('E38').Value = '  -8.58%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.09151'
$ws.Range('E39').Value = '  -6.46%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '1.452'
$ws.Range('E40').Value = '  -2.01%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '9.865'
$ws.Range('E41').Value = '  -11.47%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.7525'
$ws.Range('E42').Value = '  -9.47%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '12.95'
$ws.Range('E43').Value = '  -7.88%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '15.93'
$ws.Range('E44').Value = '  -7.17%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.6852'
$ws.Range('E45').Value = '  -7.45%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '2.456'
$ws.Range('E46').Value = '  -8.69%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '4.087'
$ws.Range('E47').Value = '  -3.94%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '1.001'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.08305'
$ws.Range('E49').Value = '  -9.59%  '
$ws.Range('B50').Value = 'Quant'
$ws.Range('C50').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '132.56'
$ws.Range('E50').Value = '  -5.65%  '
$ws.Range('B51').Value = 'Flow'
$ws.Range('C51').Value = 'https://coinranking.com/coin/QQ0NCmjVq+flow-flow'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '1.247'
$ws.Range('E51').Value = '  -13.43%  '
